# Auto-applied updates to cryptos worksheet (price/volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.147.89"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.305.40"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'316.20"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'40.27"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'0.0915"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "'0.983"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "'15.48"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").Value = "2.656.49"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "2.304.21"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "42.159.18"
$ws.Range("D19").Value = "'7.73"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").Value = "'73.19"
$ws.Range("E21").Value = "  -3.74%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'262.33"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").Value = "'9.89"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").Value = "'11.05"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").Value = "'22.92"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'36.75"
$ws.Range("E30").Value = "  +4.03%  "
$ws.Range("D31").Value = "'166.36"
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").Value = "'0.0899"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  +6.70%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +3.45%  "
$ws.Range("D38").Value = "'2.96"
$ws.Range("E38").Value = "  +14.07%  "
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "'3.64"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").Value = "'100.65"
$ws.Range("E41").Value = "  +18.42%  "
$ws.Range("D42").Value = "'1.49"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("D43").Value = "'72.07"
$ws.Range("E43").Value = "  +3.95%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'12.38"
$ws.Range("E46").Value = "  +6.04%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "'81.12"
$ws.Range("E47").Value = "  +11.89%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'114.31"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "'9.24"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("D50").Value = "'5.36"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  +3.59%  "
